$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the User ID column values (A2:A4) to shift them by +1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4

# Move selection to A4, matching the final cursor position in the diff
$ws.Range("A4").Select()
